# Refresh the "fixed" date fields that live on the slide master, every
# custom (slide) layout, and the notes master. PowerPoint stamps these the
# day the "Insert > Header & Footer > Fixed" date was (re)applied; the
# author re-applied it a week later (22/10/2018 -> 29/10/2018).
# ppPlaceholderDate = 16 identifies the date placeholder regardless of the
# auto-generated shape name (it differs per layout).
$ppPlaceholderDate = 16

# --- Slide master: date placeholder (en-US, M/D/YYYY) ---
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $shp.TextFrame.TextRange.Text = "10/29/2018"
    }
}

# --- Every custom (slide) layout: same date placeholder ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = "10/29/2018"
        }
    }
}

# --- Notes master: date placeholder (en-SG, D/M/YYYY, auto-update fld) ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $shp.TextFrame.TextRange.Text = "29/10/2018"
    }
}

# Slide 33 ("End of chapter"): drop the two leftover empty layout
# placeholders (Title 1 / Content Placeholder 2) that were never filled in;
# keep the actual "TextBox 3" / slide-number shapes.
$last = $p.Slides.Item($p.Slides.Count)
$last.Shapes.Item("Title 1").Delete()
$last.Shapes.Item("Content Placeholder 2").Delete()
